# Update cryptocurrency price/volume data to reflect the latest snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'64.955.85"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.74%  "
$c = $ws.Range("D3")
$c.Value = "'3.537.66"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.72%  "
$ws.Range("E4").Value = "  -0.37%  "
$c = $ws.Range("D5")
$c.Value = "'594.00"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.37%  "
$c = $ws.Range("D6")
$c.Value = "'136.74"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.35%  "
$c = $ws.Range("D7")
$c.Value = "'3.536.08"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.70%  "
$ws.Range("E8").Value = "  -0.12%  "
$c = $ws.Range("D9")
$c.Value = "'0.492"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +2.34%  "
$ws.Range("E10").Value = "  +3.04%  "
$c = $ws.Range("D11")
$c.Value = "'6.94"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.64%  "
$c = $ws.Range("D12")
$c.Value = "'0.383"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.05%  "
$c = $ws.Range("D13")
$c.Value = "'4.131.64"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +3.32%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D14")
$c.Value = "'0.0000181"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +2.87%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D15")
$c.Value = "'27.15"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.32%  "
$c = $ws.Range("D16")
$c.Value = "'3.530.72"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +1.33%  "
$c = $ws.Range("D18")
$c.Value = "'64.835.87"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$c = $ws.Range("D19")
$c.Value = "'10.16"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +7.45%  "
$ws.Range("E20").Value = "  +1.09%  "
$c = $ws.Range("D21")
$c.Value = "'14.23"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +5.67%  "
$c = $ws.Range("D22")
$c.Value = "'388.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("E23").Value = "  +6.14%  "
$c = $ws.Range("D24")
$c.Value = "'3.670.40"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.25%  "
$c = $ws.Range("D25")
$c.Value = "'73.92"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.87%  "
$c = $ws.Range("D26")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$c = $ws.Range("D27")
$c.Value = "'0.0000113"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +8.61%  "
$ws.Range("E28").Value = "  +7.31%  "
$ws.Range("E29").Value = "  +0.11%  "
$c = $ws.Range("D30")
$c.Value = "'2.27"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +4.76%  "
$c = $ws.Range("D31")
$c.Value = "'8.18"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.47%  "
$c = $ws.Range("D32")
$c.Value = "'3.541.36"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +3.28%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +4.17%  "
$ws.Range("E35").Value = "  +14.72%  "
$ws.Range("E36").Value = "  +2.21%  "
$c = $ws.Range("D37")
$c.Value = "'169.31"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("E38").Value = "  +8.23%  "
$ws.Range("E39").Value = "  +3.41%  "
$c = $ws.Range("D40")
$c.Value = "'4.93"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +8.71%  "
$c = $ws.Range("D41")
$c.Value = "'0.0802"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +7.22%  "
$ws.Range("E42").Value = "  +1.89%  "
$c = $ws.Range("D43")
$c.Value = "'26.20"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +16.66%  "
$c = $ws.Range("D44")
$c.Value = "'42.64"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.48%  "
$ws.Range("E45").Value = "  -0.42%  "
$c = $ws.Range("D46")
$c.Value = "'4.41"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.93%  "
$ws.Range("E47").Value = "  +7.94%  "
$ws.Range("E48").Value = "  +3.20%  "
$ws.Range("E49").Value = "  +7.03%  "
$c = $ws.Range("D50")
$c.Value = "'2.419.20"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +11.75%  "
$c = $ws.Range("D51")
$c.Value = "'302.97"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +13.14%  "
